$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'MATRICE D''IMPLICATION : PROJET PROGAV'
$ws.Range("B4").Value = "Nom `nProjet :"
$ws.Range('C4').Value = 'projet-ensemenc-bernard-adelmard'
$ws.Range('E4').Value = 'date: 23/05/2025'
$ws.Range('B5').Value = 'MEMBRES DE L''EQUIPE:'
$ws.Range('C7').Value = 'Réflexion préalable'
$ws.Range('D7').Value = 'Julien / Nathan'
$ws.Range('E7').Value = '60 / 40'
$ws.Range('C8').Value = 'MenuChoix'
$ws.Range('D8').Value = 'Nathan'
$ws.Range('C9').Value = 'Plantes et ses sous classes'
$ws.Range('D9').Value = 'Julien / Nathan'
$ws.Range('E9').Value = '50 / 50'
$ws.Range('C10').Value = 'Jeu'
$ws.Range('D10').Value = 'Julien / Nathan'
$ws.Range('E10').Value = '50 / 50'
$ws.Range('C11').Value = 'Joueur'
$ws.Range('D11').Value = 'Julien'
$ws.Range('C12').Value = 'AffichageParcelle'
$ws.Range('D12').Value = 'Nathan'
$ws.Range('C13').Value = 'AffichageChargement'
$ws.Range('D13').Value = 'Julien / Nathan'
$ws.Range('E13').Value = '30 / 70'
$ws.Range('C14').Value = 'ModeUrgence'
$ws.Range('D14').Value = 'Julien / Nathan'
$ws.Range('E14').Value = '50 / 50 '
$ws.Range('C15').Value = 'Parcelle et ses sous classes'
$ws.Range('D15').Value = 'Julien / Nathan'
$ws.Range('E15').Value = '50 / 50 '
$ws.Range('C16').Value = 'Partie'
$ws.Range('D16').Value = 'Julien / Nathan'
$ws.Range('E16').Value = '35 / 65'
$ws.Range('C17').Value = 'GestionSaison'
$ws.Range('D17').Value = 'Julien / Nathan'
$ws.Range('E17').Value = '60 /40'
$ws.Range('C18').Value = 'SauvegardeManager'
$ws.Range('D18').Value = 'Julien / Nathan'
$ws.Range('E18').Value = '70 / 30'
$ws.Range('C19').Value = 'Corrections de bugs et tests'
$ws.Range('D19').Value = 'Julien / Nathan'
$ws.Range('E19').Value = '50 / 50'
$ws.Range('C20').Value = 'Equilibrage du jeu '
$ws.Range('D20').Value = 'Julien / Nathan'
$ws.Range('E20').Value = '50 / 50'

# Re-apply the style from a data row onto the sentinel last-row cell in column E
# (mirrors a stray formatting artifact left over in the source workbook)
$ws.Cells.Item(9, 5).Copy()
$ws.Cells.Item(1048576, 5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore selection near the filled-in matrix, as left by the author
$ws.Range("C22").Select()
